$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column (H) mirroring the existing header formatting (copied from G1,
# the last existing header cell) and fill in the "1" flag values for each data row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
